# feat: add 2022-Q1 data
#
# The workbook has per-quarter detail sheets (basic "which funds hold this
# stock" breakdown) plus a trailing "总计" (Total) roll-up sheet. This change:
#   1. Inserts a new "2022-Q1" detail sheet right before "总计".
#   2. Rebuilds "总计" with a new first data row for 2022-Q1 (existing rows
#      shift down by one), so the new sheet ends up with sheetId=6 (reusing
#      the id freed by deleting/recreating 总计) and 总计 becomes sheetId=7 -
#      matching how Excel hands out sheetId's when sheets are deleted and
#      re-added in order.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 0. Grab a cell that already has the workbook's "header / index column"
#    style (bold, thin border, centered/top-aligned) so new header and
#    index cells can reuse it via copy/paste-of-formats instead of
#    re-deriving a (slightly different) style through Font/Border writes.
# ---------------------------------------------------------------------------
$styleSource = $wb.Worksheets.Item("2021-Q4").Range("B1")

# ---------------------------------------------------------------------------
# 1. Remove the old "总计" sheet, then recreate "2022-Q1" and "总计" in the
#    right order so tab order + sheetId numbering match the target.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Delete()

$lastDetailSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1Sheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastDetailSheet)
$q1Sheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1Sheet)
$totalSheet.Name = "总计"

# ---------------------------------------------------------------------------
# 2. Populate "2022-Q1" - same layout as the other per-quarter detail sheets
#    (基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
# ---------------------------------------------------------------------------
$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q1Sheet.Cells.Item(1, $col)
    $styleSource.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $q1Headers[$col - 2]
}
$excel.CutCopyMode = $false

$q1Rows = @(
    @("000603", "易方达创新驱动灵活配置混合", "71.53", "90.89", "3.40", "2.4320", 6),
    @("519692", "交银成长混合A",              "23.39", "82.72", "2.70", "0.6315", 9),
    @("960016", "交银成长混合H",              "23.39", "82.72", "2.70", "0.6315", 9),
    @("519694", "交银蓝筹混合",                "16.14", "82.45", "2.69", "0.4342", 10),
    @("519158", "新华趋势领航混合",            "4.25",  "89.82", "3.28", "0.1394", 8),
    @("730001", "方正富邦创新动力混合A",        "0.54",  "87.53", "4.33", "0.0234", 6),
    @("007046", "方正富邦创新动力混合C",        "0.33",  "87.53", "4.33", "0.0143", 6)
)

for ($i = 0; $i -lt $q1Rows.Count; $i++) {
    $r = $i + 2
    $row = $q1Rows[$i]

    $idxCell = $q1Sheet.Cells.Item($r, 1)
    $styleSource.Copy()
    $idxCell.PasteSpecial(-4122)
    $idxCell.Value = $i

    $q1Sheet.Cells.Item($r, 2).Value = "'" + $row[0]   # 基金代码 (text, keep leading zeros)
    $q1Sheet.Cells.Item($r, 3).Value = $row[1]          # 基金名称
    $q1Sheet.Cells.Item($r, 4).Value = "'" + $row[2]    # 基金规模 (text)
    $q1Sheet.Cells.Item($r, 5).Value = "'" + $row[3]    # 股票总仓位 (text)
    $q1Sheet.Cells.Item($r, 6).Value = "'" + $row[4]    # 仓位占比 (text)
    $q1Sheet.Cells.Item($r, 7).Value = "'" + $row[5]    # 持有市值(亿元) (text)
    $q1Sheet.Cells.Item($r, 8).Value = $row[6]          # 仓位排名 (number)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Rebuild "总计" - same 日期/持有数量(只)/持有市值(亿元) roll-up, now with
#    2022-Q1 inserted as the newest (first) row.
# ---------------------------------------------------------------------------
$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($col = 2; $col -le 4; $col++) {
    $cell = $totalSheet.Cells.Item(1, $col)
    $styleSource.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $totalHeaders[$col - 2]
}
$excel.CutCopyMode = $false

$totalRows = @(
    @("2022-Q1", 7,  4.31),
    @("2021-Q4", 24, 15.47),
    @("2021-Q3", 31, 18.7),
    @("2021-Q2", 19, 9.32),
    @("2021-Q1", 10, 4.37),
    @("2020-Q4", 17, 5.51)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]

    $idxCell = $totalSheet.Cells.Item($r, 1)
    $styleSource.Copy()
    $idxCell.PasteSpecial(-4122)
    $idxCell.Value = $i

    $totalSheet.Cells.Item($r, 2).Value = $row[0]  # 日期 (plain text, not numeric-looking)
    $totalSheet.Cells.Item($r, 3).Value = $row[1]  # 持有数量(只) (number)
    $totalSheet.Cells.Item($r, 4).Value = $row[2]  # 持有市值(亿元) (number)
}
$excel.CutCopyMode = $false

$q1Sheet.Range("A1").Select()
$totalSheet.Range("A1").Select()
